$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44309
$ws.Range("M2").Value2 = 40

# Row 3
$ws.Range("D3").Value2 = 44309
$ws.Range("M3").Value2 = 70

# Row 4
$ws.Range("D4").Value2 = 44305
$ws.Range("M4").Value2 = 50

# Row 5
$ws.Range("D5").Value2 = 44305
$ws.Range("M5").Value2 = 60
